$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 54

$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-10-08"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 52.77000045776367
$ws.Cells.Item($newRow, 3).Value = 681.5499877929688
$ws.Cells.Item($newRow, 4).Value = 341.6499938964844
